$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S2").Value = 56.53109005
$ws.Range("T2").Value = 68.72117165
$ws.Range("U2").Value = 83.5398615
$ws.Range("V2").Value = 101.5539795
$ws.Range("W2").Value = 123.4525719
$ws.Range("X2").Value = 150.0732673
$ws.Range("Y2").Value = 182.4343162
$ws.Range("Z2").Value = 199.2018542
$ws.Range("AA2").Value = 217.5104967
$ws.Range("AB2").Value = 237.5018864
$ws.Range("AC2").Value = 259.3306848
$ws.Range("AD2").Value = 283.165768
$ws.Range("AE2").Value = 321.6038221
$ws.Range("AF2").Value = 365.2596113
$ws.Range("AG2").Value = 414.841412
$ws.Range("AH2").Value = 471.1536446
$ws.Range("AI2").Value = 535.1099249
$ws.Range("AJ2").Value = 584.2918781
$ws.Range("AK2").Value = 637.9941446
$ws.Range("AL2").Value = 696.6321865
$ws.Range("AM2").Value = 760.6596507
$ws.Range("AN2").Value = 830.5718792
$ws.Range("AO2").Value = 912.1946695
$ws.Range("AP2").Value = 1001.838776
$ws.Range("AQ2").Value = 1100.292478
$ws.Range("AR2").Value = 1208.421521
$ws.Range("AS2").Value = 1327.176729
